$d = $word.ActiveDocument

# The document contains a single results table. Rows 2-13 (1-indexed, i.e.
# skipping the header row) hold the "% Variance" values in column 5. We
# replace each one with its debugged value, matched by row/column position
# so duplicate percentage strings (e.g. "2.1% ") are handled unambiguously.

$table = $d.Tables.Item(1)

$changes = @(
    @{ Row = 2;  Col = 5; Old = "2.9% "; New = "6.9% " },
    @{ Row = 3;  Col = 5; Old = "1.8% "; New = "4.2% " },
    @{ Row = 4;  Col = 5; Old = "7.7% "; New = "18.5%" },
    @{ Row = 5;  Col = 5; Old = "29.3%"; New = "70.4%" },
    @{ Row = 6;  Col = 5; Old = "3.0% "; New = "8.2% " },
    @{ Row = 7;  Col = 5; Old = "2.1% "; New = "5.7% " },
    @{ Row = 8;  Col = 5; Old = "10.1%"; New = "27.5%" },
    @{ Row = 9;  Col = 5; Old = "21.4%"; New = "58.5%" },
    @{ Row = 10; Col = 5; Old = "2.1% "; New = "9.5% " },
    @{ Row = 11; Col = 5; Old = "1.7% "; New = "7.7% " },
    @{ Row = 12; Col = 5; Old = "4.7% "; New = "21.8%" },
    @{ Row = 13; Col = 5; Old = "13.3%"; New = "61.1%" }
)

foreach ($ch in $changes) {
    $cell = $table.Cell($ch.Row, $ch.Col)
    $rng = $cell.Range
    # A table-cell Range includes the trailing end-of-cell mark (and
    # paragraph mark); trim it off (wdCharacter = 1) before comparing with /
    # overwriting the visible text, so we don't clobber the cell structure.
    [void]$rng.MoveEnd(1, -1)
    if ($rng.Text -ne $ch.Old) {
        throw "Unexpected text in row $($ch.Row): [$($rng.Text)] (expected [$($ch.Old)])"
    }
    $rng.Text = $ch.New
}
